$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Range("AG100")

$scratch.Formula = '="320018511882"'
$scratch.Copy()
$ws.Range("C2").PasteSpecial(-4163)

$scratch.Formula = '="320018511996"'
$scratch.Copy()
$ws.Range("C3").PasteSpecial(-4163)

$scratch.Formula = '="320018512021"'
$scratch.Copy()
$ws.Range("C4").PasteSpecial(-4163)

$scratch.Formula = '="320018512098"'
$scratch.Copy()
$ws.Range("C5").PasteSpecial(-4163)
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)

$scratch.Formula = '="320018512135"'
$scratch.Copy()
$ws.Range("C6").PasteSpecial(-4163)
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)

$scratch.Formula = '="320018512205"'
$scratch.Copy()
$ws.Range("C7").PasteSpecial(-4163)
$scratch.Copy()
$ws.Range("D7").PasteSpecial(-4163)

$scratch.Formula = '="320018512238"'
$scratch.Copy()
$ws.Range("C8").PasteSpecial(-4163)

$scratch.Formula = '="320018512250"'
$scratch.Copy()
$ws.Range("C9").PasteSpecial(-4163)

$scratch.Formula = '="320018512282"'
$scratch.Copy()
$ws.Range("C10").PasteSpecial(-4163)

$scratch.Formula = '="320018512319"'
$scratch.Copy()
$ws.Range("C11").PasteSpecial(-4163)

$scratch.Formula = '="320018512400"'
$scratch.Copy()
$ws.Range("C12").PasteSpecial(-4163)

$scratch.Formula = '="320018503130"'
$scratch.Copy()
$ws.Range("C13").PasteSpecial(-4163)
$scratch.Copy()
$ws.Range("D13").PasteSpecial(-4163)

$scratch.Formula = '="320018503163"'
$scratch.Copy()
$ws.Range("C14").PasteSpecial(-4163)
$scratch.Copy()
$ws.Range("D14").PasteSpecial(-4163)

$scratch.Formula = '="320018503185"'
$scratch.Copy()
$ws.Range("C15").PasteSpecial(-4163)
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)

$scratch.Formula = '="320018503211"'
$scratch.Copy()
$ws.Range("C16").PasteSpecial(-4163)
$scratch.Copy()
$ws.Range("D16").PasteSpecial(-4163)

$scratch.Formula = '="320018503244"'
$scratch.Copy()
$ws.Range("C17").PasteSpecial(-4163)
$scratch.Copy()
$ws.Range("D17").PasteSpecial(-4163)

$scratch.Formula = '="320018503288"'
$scratch.Copy()
$ws.Range("C18").PasteSpecial(-4163)

$scratch.Formula = '="320018503303"'
$scratch.Copy()
$ws.Range("C19").PasteSpecial(-4163)

$scratch.Formula = '="320018503336"'
$scratch.Copy()
$ws.Range("C20").PasteSpecial(-4163)

$scratch.Formula = '="320018503358"'
$scratch.Copy()
$ws.Range("C21").PasteSpecial(-4163)

$scratch.Formula = '="320018503380"'
$scratch.Copy()
$ws.Range("C22").PasteSpecial(-4163)

$scratch.ClearContents()
$excel.CutCopyMode = $false
